$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 316.625
$ws.Range("I2").Value = 207.2
$ws.Range("J2").Value = 499
$ws.Range("K2").Value = 207.2
$ws.Range("L2").Value = 499
$ws.Range("M2").Value = -94.19999999999999
$ws.Range("N2").Value = -725
$ws.Range("H5").Value = 1000619.4
$ws.Range("I5").Value = 1111782.6
$ws.Range("K5").Value = 1111782.6
$ws.Range("M5").Value = -1111667.6
$ws.Range("H17").Value = 5475.0605
$ws.Range("J17").Value = 5475.0605
$ws.Range("L17").Value = 16425.1815
$ws.Range("N17").Value = -16761.1815
$ws.Range("H29").Value = 314
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H40").Value = 5750
$ws.Range("I40").Value = 3699.6667
$ws.Range("J40").Value = 6775.1665
$ws.Range("K40").Value = 3699.6667
$ws.Range("L40").Value = 6775.1665
$ws.Range("M40").Value = -3524.6667
$ws.Range("N40").Value = -7125.1665
$ws.Range("H41").Value = 395.7143
$ws.Range("I41").Value = 426.83334
$ws.Range("K41").Value = 426.83334
$ws.Range("M41").Value = 13.16665999999998
$ws.Range("H58").Value = 429.9
$ws.Range("I58").Value = 366.66666
$ws.Range("J58").Value = 999
$ws.Range("K58").Value = 1099.99998
$ws.Range("L58").Value = 2997
$ws.Range("M58").Value = -949.9999800000001
$ws.Range("N58").Value = -3297
$ws.Range("H61").Value = 396
$ws.Range("I61").Value = 396
$ws.Range("K61").Value = 1188
$ws.Range("M61").Value = -1016
$ws.Range("H92").Value = 3322.3125
$ws.Range("I92").Value = 2928.889
$ws.Range("K92").Value = 2928.889
$ws.Range("M92").Value = -1680.889
$ws.Range("H98").Value = 344089.28
$ws.Range("I98").Value = 2509.2727
$ws.Range("K98").Value = 2509.2727
$ws.Range("M98").Value = -1011.2727
$ws.Range("H112").Value = 2062.5557
$ws.Range("I112").Value = 1500
$ws.Range("J112").Value = 2132.875
$ws.Range("K112").Value = 4500
$ws.Range("L112").Value = 6398.625
$ws.Range("M112").Value = -3392
$ws.Range("N112").Value = -8614.625
$ws.Range("H122").Value = 344089.28
$ws.Range("I122").Value = 2509.2727
$ws.Range("K122").Value = 7527.8181
$ws.Range("M122").Value = -5077.8181
$ws.Range("H132").Value = 1742.1666
$ws.Range("I132").Value = 1828.2667
$ws.Range("J132").Value = 1311.6666
$ws.Range("K132").Value = 5484.800099999999
$ws.Range("L132").Value = 3934.9998
$ws.Range("M132").Value = -2954.800099999999
$ws.Range("N132").Value = -8994.9998
$ws.Range("H137").Value = 17860146
$ws.Range("I137").Value = 58826084
$ws.Range("J137").Value = 3197.8718
$ws.Range("K137").Value = 176478252
$ws.Range("L137").Value = 9593.615399999999
$ws.Range("M137").Value = -176475702
$ws.Range("N137").Value = -14693.6154
$ws.Range("H138").Value = 3547.453
$ws.Range("I138").Value = 2097.8333
$ws.Range("J138").Value = 4292.971
$ws.Range("K138").Value = 6293.499899999999
$ws.Range("L138").Value = 12878.913
$ws.Range("M138").Value = -1153.499899999999
$ws.Range("N138").Value = -23158.913
$ws.Range("H141").Value = 1821.9
$ws.Range("I141").Value = 1422.4667
$ws.Range("K141").Value = 4267.4001
$ws.Range("M141").Value = 912.5999000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3473.2097
$ws.Range("I32").Value = 3202.2788
$ws.Range("K32").Value = 3202.2788
$ws.Range("M32").Value = -2915.2788
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()
$ws.Range("H74").Value = 7248995
$ws.Range("I74").Value = 9805830
$ws.Range("J74").Value = 4630.25
$ws.Range("K74").Value = 9805830
$ws.Range("L74").Value = 4630.25
$ws.Range("M74").Value = -9804956
$ws.Range("N74").Value = -6378.25
$ws.Range("H77").Value = 7248995
$ws.Range("I77").Value = 9805830
$ws.Range("J77").Value = 4630.25
$ws.Range("K77").Value = 49029150
$ws.Range("L77").Value = 23151.25
$ws.Range("M77").Value = -49024782
$ws.Range("N77").Value = -31887.25
$ws.Range("H102").Value = 2231.2632
$ws.Range("I102").Value = 1618.375
$ws.Range("K102").Value = 1618.375
$ws.Range("M102").Value = 3.625
$ws.Range("H132").Value = 6281.486
$ws.Range("I132").Value = 4669.76
$ws.Range("K132").Value = 14009.28
$ws.Range("M132").Value = -11479.28

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3482.152
$ws.Range("I20").Value = 2792.6428
$ws.Range("J20").Value = 4554.722
$ws.Range("K20").Value = 2792.6428
$ws.Range("L20").Value = 4554.722
$ws.Range("M20").Value = -2545.6428
$ws.Range("N20").Value = -5048.722
$ws.Range("H86").Value = 3212.4443
$ws.Range("I86").Value = 2151.4285
$ws.Range("J86").Value = 6926
$ws.Range("K86").Value = 2151.4285
$ws.Range("L86").Value = 6926
$ws.Range("M86").Value = -1028.4285
$ws.Range("N86").Value = -9172
$ws.Range("H89").Value = 3212.4443
$ws.Range("I89").Value = 2151.4285
$ws.Range("J89").Value = 6926
$ws.Range("K89").Value = 10757.1425
$ws.Range("L89").Value = 34630
$ws.Range("M89").Value = -5141.1425
$ws.Range("N89").Value = -45862
$ws.Range("H99").Value = 2086.3333
$ws.Range("I99").Value = 1797
$ws.Range("K99").Value = 1797
$ws.Range("M99").Value = -299
$ws.Range("H103").Value = 2750
$ws.Range("J103").Value = 2750
$ws.Range("L103").Value = 2750
$ws.Range("N103").Value = -5094
$ws.Range("H105").Value = 9958.031999999999
$ws.Range("I105").Value = 10811.6
$ws.Range("K105").Value = 10811.6
$ws.Range("M105").Value = -9064.6
$ws.Range("H107").Value = 1619.8636
$ws.Range("I107").Value = 1268.6428
$ws.Range("K107").Value = 1268.6428
$ws.Range("M107").Value = 651.3571999999999
$ws.Range("H134").Value = 1170.8167
$ws.Range("I134").Value = 756.25
$ws.Range("K134").Value = 2268.75
$ws.Range("M134").Value = 266.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1129.1578
$ws.Range("I7").Value = 2323.1428
$ws.Range("K7").Value = 2323.1428
$ws.Range("M7").Value = -2210.1428
$ws.Range("H22").Value = 1033.2
$ws.Range("I22").Value = 318.8889
$ws.Range("J22").Value = 2104.6667
$ws.Range("K22").Value = 318.8889
$ws.Range("L22").Value = 2104.6667
$ws.Range("M22").Value = 31.11110000000002
$ws.Range("N22").Value = -2804.6667
$ws.Range("H31").Value = 53452.316
$ws.Range("I31").Value = 4611.273
$ws.Range("J31").Value = 102293.37
$ws.Range("K31").Value = 4611.273
$ws.Range("L31").Value = 102293.37
$ws.Range("M31").Value = -4316.273
$ws.Range("N31").Value = -102883.37
$ws.Range("H34").Value = 53452.316
$ws.Range("I34").Value = 4611.273
$ws.Range("J34").Value = 102293.37
$ws.Range("K34").Value = 4611.273
$ws.Range("L34").Value = 102293.37
$ws.Range("M34").Value = -4409.273
$ws.Range("N34").Value = -102697.37
$ws.Range("H55").Value = 25000
$ws.Range("J55").Value = 25000
$ws.Range("L55").Value = 25000
$ws.Range("N55").Value = -25630
$ws.Range("H58").Value = 3230.075
$ws.Range("I58").Value = 2028.9667
$ws.Range("K58").Value = 2028.9667
$ws.Range("M58").Value = -1825.9667
$ws.Range("H105").Value = 4573.857
$ws.Range("I105").Value = 1998.75
$ws.Range("J105").Value = 8007.3335
$ws.Range("K105").Value = 1998.75
$ws.Range("L105").Value = 8007.3335
$ws.Range("M105").Value = -251.75
$ws.Range("N105").Value = -11501.3335
$ws.Range("H107").Value = 1194.5769
$ws.Range("I107").Value = 1075.9375
$ws.Range("K107").Value = 1075.9375
$ws.Range("M107").Value = 844.0625
$ws.Range("H109").Value = 58626.168
$ws.Range("I109").Value = 45269
$ws.Range("J109").Value = 61297.6
$ws.Range("K109").Value = 45269
$ws.Range("L109").Value = 61297.6
$ws.Range("M109").Value = -44229
$ws.Range("N109").Value = -63377.6
$ws.Range("H122").Value = 5145.5
$ws.Range("I122").Value = 2446.6
$ws.Range("K122").Value = 7339.799999999999
$ws.Range("M122").Value = -4889.799999999999
$ws.Range("H132").Value = 2911.9707
$ws.Range("I132").Value = 2404.9678
$ws.Range("K132").Value = 7214.903399999999
$ws.Range("M132").Value = -4684.903399999999
$ws.Range("H134").Value = 4971.415
$ws.Range("I134").Value = 4295.3696
$ws.Range("J134").Value = 9414
$ws.Range("K134").Value = 12886.1088
$ws.Range("L134").Value = 28242
$ws.Range("M134").Value = -10351.1088
$ws.Range("N134").Value = -33312
$ws.Range("H136").Value = 3230.075
$ws.Range("I136").Value = 2028.9667
$ws.Range("K136").Value = 6086.9001
$ws.Range("M136").Value = -3536.9001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 7683222.5
$ws.Range("I4").Value = 6322153
$ws.Range("K4").Value = 18966459
$ws.Range("M4").Value = -18966347
$ws.Range("H23").Value = 425.5
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 425.5
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 1276.5
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -1746.5
$ws.Range("H63").Value = 13754.846
$ws.Range("I63").Value = 3265.3333
$ws.Range("J63").Value = 16901.7
$ws.Range("K63").Value = 9795.999899999999
$ws.Range("L63").Value = 50705.10000000001
$ws.Range("M63").Value = -9046.999899999999
$ws.Range("N63").Value = -52203.10000000001
$ws.Range("H66").Value = 13754.846
$ws.Range("I66").Value = 3265.3333
$ws.Range("J66").Value = 16901.7
$ws.Range("K66").Value = 29387.9997
$ws.Range("L66").Value = 152115.3
$ws.Range("M66").Value = -25643.9997
$ws.Range("N66").Value = -159603.3
$ws.Range("H80").Value = 4001
$ws.Range("J80").Value = 4334.6665
$ws.Range("L80").Value = 13003.9995
$ws.Range("N80").Value = -14875.9995
$ws.Range("H83").Value = 4001
$ws.Range("J83").Value = 4334.6665
$ws.Range("L83").Value = 39011.9985
$ws.Range("N83").Value = -48371.9985
$ws.Range("H86").Value = 2939.2222
$ws.Range("J86").Value = 3637.1428
$ws.Range("L86").Value = 10911.4284
$ws.Range("N86").Value = -13283.4284
$ws.Range("H87").Value = 8887.666999999999
$ws.Range("I87").Value = 8887.666999999999
$ws.Range("K87").Value = 26663.001
$ws.Range("M87").Value = -25415.001
$ws.Range("H89").Value = 2939.2222
$ws.Range("J89").Value = 3637.1428
$ws.Range("L89").Value = 32734.2852
$ws.Range("N89").Value = -44590.2852
$ws.Range("H90").Value = 8887.666999999999
$ws.Range("I90").Value = 8887.666999999999
$ws.Range("K90").Value = 79989.003
$ws.Range("M90").Value = -73749.003
$ws.Range("H103").Value = 2259.8667
$ws.Range("I103").Value = 2850.8572
$ws.Range("J103").Value = 1742.75
$ws.Range("K103").Value = 8552.571599999999
$ws.Range("L103").Value = 5228.25
$ws.Range("M103").Value = -7673.571599999999
$ws.Range("N103").Value = -6986.25
$ws.Range("H108").Value = 2139.7
$ws.Range("I108").Value = 263
$ws.Range("K108").Value = 789
$ws.Range("M108").Value = 2091
$ws.Range("H121").Value = 1141.6428
$ws.Range("I121").Value = 1347.8572
$ws.Range("J121").Value = 935.4286
$ws.Range("K121").Value = 4043.5716
$ws.Range("L121").Value = 2806.2858
$ws.Range("M121").Value = -2733.5716
$ws.Range("N121").Value = -5426.2858
$ws.Range("H129").Value = 3971485.8
$ws.Range("I129").Value = 722.8182
$ws.Range("J129").Value = 8339325
$ws.Range("K129").Value = 2168.4546
$ws.Range("L129").Value = 25017975
$ws.Range("M129").Value = 2831.5454
$ws.Range("N129").Value = -25027975
$ws.Range("H141").Value = 5684.0967
$ws.Range("I141").Value = 3361.2727
$ws.Range("J141").Value = 11362.111
$ws.Range("K141").Value = 10083.8181
$ws.Range("L141").Value = 34086.333
$ws.Range("M141").Value = -4903.8181
$ws.Range("N141").Value = -44446.333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7102.143
$ws.Range("I80").Value = 2401.3333
$ws.Range("K80").Value = 2401.3333
$ws.Range("M80").Value = -1403.3333
$ws.Range("H83").Value = 7102.143
$ws.Range("I83").Value = 2401.3333
$ws.Range("K83").Value = 12006.6665
$ws.Range("M83").Value = -7014.666499999999
$ws.Range("H97").Value = 738.0833
$ws.Range("I97").Value = 783.8095
$ws.Range("K97").Value = 783.8095
$ws.Range("M97").Value = -287.8095
$ws.Range("H101").Value = 49499.5
$ws.Range("J101").Value = 49499.5
$ws.Range("L101").Value = 49499.5
$ws.Range("N101").Value = -55989.5
$ws.Range("H122").Value = 5203.2354
$ws.Range("I122").Value = 5075.5815
$ws.Range("K122").Value = 15226.7445
$ws.Range("M122").Value = -12776.7445
$ws.Range("H126").Value = 3579.2354
$ws.Range("I126").Value = 2712.182
$ws.Range("J126").Value = 5168.8335
$ws.Range("K126").Value = 8136.545999999999
$ws.Range("L126").Value = 15506.5005
$ws.Range("M126").Value = -5666.545999999999
$ws.Range("N126").Value = -20446.5005
$ws.Range("H129").Value = 60000
$ws.Range("J129").Value = 60000
$ws.Range("L129").Value = 60000
$ws.Range("N129").Value = -70000
$ws.Range("H132").Value = 5981.636
$ws.Range("I132").Value = 6818.8
$ws.Range("K132").Value = 20456.4
$ws.Range("M132").Value = -17926.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 991.4583
$ws.Range("I16").Value = 991.4583
$ws.Range("K16").Value = 991.4583
$ws.Range("M16").Value = -821.4583
$ws.Range("H40").Value = 6076.4
$ws.Range("I40").Value = 5783.48
$ws.Range("K40").Value = 5783.48
$ws.Range("M40").Value = -5647.48
$ws.Range("H55").Value = 3574170.8
$ws.Range("I55").Value = 7143655.5
$ws.Range("J55").Value = 4685.857
$ws.Range("K55").Value = 7143655.5
$ws.Range("L55").Value = 4685.857
$ws.Range("M55").Value = -7143482.5
$ws.Range("N55").Value = -5031.857
$ws.Range("H132").Value = 3358.5173
$ws.Range("I132").Value = 2956.0637
$ws.Range("K132").Value = 8868.1911
$ws.Range("M132").Value = -6338.1911
$ws.Range("H136").Value = 3767.2239
$ws.Range("I136").Value = 2886.6965
$ws.Range("J136").Value = 8249.909
$ws.Range("K136").Value = 8660.0895
$ws.Range("L136").Value = 24749.727
$ws.Range("M136").Value = -6110.0895
$ws.Range("N136").Value = -29849.727

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 20000
$ws.Range("J54").Value = 22500
$ws.Range("L54").Value = 22500
$ws.Range("N54").Value = -23540
$ws.Range("H75").Value = 36118
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 36118
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H107").Value = 1295.4615
$ws.Range("I107").Value = 1167.3636
$ws.Range("K107").Value = 3502.0908
$ws.Range("M107").Value = -1582.0908
$ws.Range("H113").Value = 608.44446
$ws.Range("I113").Value = 405
$ws.Range("J113").Value = 811.8889
$ws.Range("K113").Value = 1215
$ws.Range("L113").Value = 2435.6667
$ws.Range("M113").Value = 955
$ws.Range("N113").Value = -6775.6667
$ws.Range("H122").Value = 3983.1924
$ws.Range("I122").Value = 4228.5
$ws.Range("K122").Value = 12685.5
$ws.Range("M122").Value = -10235.5
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H132").Value = 2226.9111
$ws.Range("I132").Value = 1218.1471
$ws.Range("K132").Value = 3654.4413
$ws.Range("M132").Value = -1124.4413
$ws.Range("H136").Value = 2140.7925
$ws.Range("I136").Value = 1611.0851
$ws.Range("K136").Value = 4833.2553
$ws.Range("M136").Value = -2283.2553
$ws.Range("H141").Value = 92669.53
$ws.Range("J141").Value = 94742.36
$ws.Range("L141").Value = 94742.36
$ws.Range("N141").Value = -105102.36
